$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.885541200637817
$ws.Range("B1").Value = 2.975127458572388
$ws.Range("C1").Value = 2.063964128494263
$ws.Range("D1").Value = 1.849264740943909
$ws.Range("E1").Value = 1.746202826499939
